$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91 (pushes old rows 91..182 down to 92..183)
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with its data
$ws.Range("A91").Value = 8
$ws.Range("B91").Value = "Terminal La Palmera de La Serena"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44586
$ws.Range("D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = 4
$ws.Range("F91").Value = 100112021
$ws.Range("G91").Value = "Ají"
$ws.Range("H91").Value = "Americana (o)"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 560
$ws.Range("K91").Value = 12000
$ws.Range("L91").Value = 13000
$ws.Range("M91").Value = 12500
$ws.Range("N91").Value = "$/caja 15 kilos"
$ws.Range("O91").Value = "Provincia de Limarí"
$ws.Range("P91").Value = 833
$ws.Range("Q91").Value = 15
$ws.Range("R91").Value = "Hortaliza"
